# This script applies the changes described by the commit
# "remove final occurrences of the word 'prio'" to the bitstream
# documentation worksheet.
#
# The underlying data rename is: a handful of VHDL/function names that
# still referred to the old "priority" (prio) signal naming scheme are
# renamed to the newer "vert"/"bus[n]" naming scheme that the rest of
# the sheet already uses (e.g. xp_bus[3]_prio[0] -> xp_bus[3]_vert[4],
# en_prio_east[1] -> en_bus_east[5], xp_prio_ns[1]_ew[1] -> xp_ns[5]_ew[5],
# and the two "priority[x]" wording fragments inside comment cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bitstream")

# --- Section CBv (rows 4-7): Function Name column (F) ---
$ws.Range("F4").Value = "xp_bus[3]_vert[4]"
$ws.Range("F5").Value = "xp_bus[2]_vert[5]"
$ws.Range("F6").Value = "xp_bus[1]_vert[4]"
$ws.Range("F7").Value = "xp_bus[0]_vert[5]"

# --- Section CBh (rows 96 & 104): Comment column (G) ---
$ws.Range("G96").Value = "enables cin sourcing from bus[5] for own adder"
$ws.Range("G104").Value = "enables cout from LUT on bus[4] as cin for other adders"

# --- Section SW (rows 108-111): Function Name column (F) now mirrors
#     the Comment column (G), using the same "xp_ns[..]_ew[..]" naming,
#     and picks up the left-aligned "Consolas" style already used by G ---
$ws.Range("F108").Value = "xp_ns[5]_ew[5]"
$ws.Range("F109").Value = "xp_ns[4]_ew[5]"
$ws.Range("F110").Value = "xp_ns[5]_ew[4]"
$ws.Range("F111").Value = "xp_ns[4]_ew[4]"

$ws.Range("F108").HorizontalAlignment = -4131
$ws.Range("F109").HorizontalAlignment = -4131
$ws.Range("F110").HorizontalAlignment = -4131
$ws.Range("F111").HorizontalAlignment = -4131

# --- Section SW (rows 132-139): Function Name column (F) now mirrors
#     the Comment column (G) naming as well ---
$ws.Range("F132").Value = "en_bus_east[5]"
$ws.Range("F133").Value = "en_bus_east[4]"
$ws.Range("F134").Value = "en_bus_west[5]"
$ws.Range("F135").Value = "en_bus_west[4]"
$ws.Range("F136").Value = "en_bus_north[5]"
$ws.Range("F137").Value = "en_bus_north[4]"
$ws.Range("F138").Value = "en_bus_south[5]"
$ws.Range("F139").Value = "en_bus_south[4]"

# --- Update the view/selection to match where the author ended up ---
$ws.Range("A90").Select()
$excel.ActiveWindow.ScrollRow = 90
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G102").Select()
